$d = $word.ActiveDocument

# The last paragraph currently holds: bookmarkStart(_GoBack) + drawing run + bookmarkEnd(_GoBack)
$lastPara = $d.Paragraphs.Last
$endOfDrawingPara = $lastPara.Range

Write-Output $endOfDrawingPara.Text
